# Add a new "EarningsReport" worksheet after the existing "IncomeReport"
# sheet and populate it with EPS / earnings-related data, mirroring the
# layout/style already used on the IncomeReport sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet after the existing (last) sheet, then rename it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "EarningsReport"

# Header row (row 1).
$ws2.Range("A1").Value = "EPS Estimate"
$ws2.Range("B1").Value = "EPS Actual"
$ws2.Range("C1").Value = "EPS Surprise"
$ws2.Range("D1").Value = "EPS Growth Quarter and Year Forecast"
$ws2.Range("E1").Value = "Revenue Growth Quarter and Year Forecast"
$ws2.Range("F1").Value = "Price Delta"

# Data row (row 2). A2/B2 look numeric but are stored as text in the
# target workbook, so force a text number format before assigning them.
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "0.51"
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "0.57"
$ws2.Range("C2").Value = 11
$ws2.Range("D2").Value = "55900, 402"
$ws2.Range("E2").Value = "35, 35"
$ws2.Range("F2").Value = 47

# Reuse the same cell formatting (fill/border/alignment) already used on
# the IncomeReport sheet's header and data rows.
$ws1.Range("A1:F1").Copy()
$ws2.Range("A1:F1").PasteSpecial(-4122)
$ws1.Range("A2:F2").Copy()
$ws2.Range("A2:F2").PasteSpecial(-4122)

# Column widths, matching the target layout as closely as this runtime's
# character-based ColumnWidth unit allows.
$ws2.Columns.Item(1).ColumnWidth = 11.780170478820798
$ws2.Columns.Item(2).ColumnWidth = 9.780298423767091
$ws2.Columns.Item(3).ColumnWidth = 11.780170478820798
$ws2.Columns.Item(4).ColumnWidth = 35.78016082763673
$ws2.Columns.Item(5).ColumnWidth = 39.78021007537841
$ws2.Columns.Item(6).ColumnWidth = 10.780234451293946

# Row heights.
$ws2.Rows.Item(1).RowHeight = 20
$ws2.Rows.Item(2).RowHeight = 30
